$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 467.73685
$ws.Range("I19").Value = 513.1667
$ws.Range("K19").Value = 513.1667
$ws.Range("M19").Value = -338.1667

$ws.Range("H33").Value = 2666
$ws.Range("I33").Value = 2700.2
$ws.Range("J33").Value = 2495
$ws.Range("K33").Value = 2700.2
$ws.Range("L33").Value = 2495
$ws.Range("M33").Value = -2471.2
$ws.Range("N33").Value = -2953

$ws.Range("H74").Value = 7132
$ws.Range("I74").Value = 4397.8184
$ws.Range("J74").Value = 11428.571
$ws.Range("K74").Value = 4397.8184
$ws.Range("L74").Value = 11428.571
$ws.Range("M74").Value = -3461.8184
$ws.Range("N74").Value = -13300.571

$ws.Range("H77").Value = 7132
$ws.Range("I77").Value = 4397.8184
$ws.Range("J77").Value = 11428.571
$ws.Range("K77").Value = 21989.092
$ws.Range("L77").Value = 57142.855
$ws.Range("M77").Value = -17309.092
$ws.Range("N77").Value = -66502.855

$ws.Range("H101").Value = 3936.6667
$ws.Range("I101").Value = 905
$ws.Range("K101").Value = 2715
$ws.Range("M101").Value = -1093

$ws.Range("H116").Value = 3070.6875
$ws.Range("I116").Value = 3070.889
$ws.Range("J116").Value = 3070.4285
$ws.Range("K116").Value = 3070.889
$ws.Range("L116").Value = 3070.4285
$ws.Range("M116").Value = 371.1109999999999
$ws.Range("N116").Value = -9954.4285

$ws.Range("H132").Value = 14539.514
$ws.Range("I132").Value = 1044.303
$ws.Range("J132").Value = 125875
$ws.Range("K132").Value = 3132.909000000001
$ws.Range("L132").Value = 377625
$ws.Range("M132").Value = -602.9090000000006
$ws.Range("N132").Value = -382685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1415.9642
$ws.Range("I74").Value = 1445.92
$ws.Range("J74").Value = 1166.3334
$ws.Range("K74").Value = 1445.92
$ws.Range("L74").Value = 1166.3334
$ws.Range("M74").Value = -571.9200000000001
$ws.Range("N74").Value = -2914.3334

$ws.Range("H77").Value = 1415.9642
$ws.Range("I77").Value = 1445.92
$ws.Range("J77").Value = 1166.3334
$ws.Range("K77").Value = 7229.6
$ws.Range("L77").Value = 5831.666999999999
$ws.Range("M77").Value = -2861.6
$ws.Range("N77").Value = -14567.667

$ws.Range("H97").Value = 1412.7778
$ws.Range("I97").Value = 1412.7778
$ws.Range("K97").Value = 1412.7778
$ws.Range("M97").Value = -916.7778000000001

$ws.Range("H110").Value = 2673.2778
$ws.Range("I110").Value = 2673.2778
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2673.2778
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -628.2777999999998
$ws.Range("N110").ClearContents()

$ws.Range("H134").Value = 59997.5
$ws.Range("J134").Value = 59997.5
$ws.Range("L134").Value = 59997.5
$ws.Range("N134").Value = -70137.5

$ws.Range("H135").Value = 73730.766
$ws.Range("J135").Value = 73730.766
$ws.Range("L135").Value = 73730.766
$ws.Range("N135").Value = -83870.766

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 274.53333
$ws.Range("J80").Value = 279.85715
$ws.Range("L80").Value = 279.85715
$ws.Range("N80").Value = -2275.85715

$ws.Range("H83").Value = 274.53333
$ws.Range("J83").Value = 279.85715
$ws.Range("L83").Value = 1399.28575
$ws.Range("N83").Value = -11383.28575

$ws.Range("H86").Value = 2440.6667
$ws.Range("I86").Value = 1547.2667
$ws.Range("J86").Value = 3557.4167
$ws.Range("K86").Value = 1547.2667
$ws.Range("L86").Value = 3557.4167
$ws.Range("M86").Value = -424.2666999999999
$ws.Range("N86").Value = -5803.4167

$ws.Range("H89").Value = 2440.6667
$ws.Range("I89").Value = 1547.2667
$ws.Range("J89").Value = 3557.4167
$ws.Range("K89").Value = 7736.3335
$ws.Range("L89").Value = 17787.0835
$ws.Range("M89").Value = -2120.3335
$ws.Range("N89").Value = -29019.0835

$ws.Range("H134").Value = 3563.724
$ws.Range("I134").Value = 1812.3158
$ws.Range("K134").Value = 5436.9474
$ws.Range("M134").Value = -2901.9474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 582.5
$ws.Range("I22").Value = 582.5
$ws.Range("K22").Value = 582.5
$ws.Range("M22").Value = -232.5

$ws.Range("H26").Value = 5999.5
$ws.Range("J26").Value = 5999.5
$ws.Range("L26").Value = 5999.5
$ws.Range("N26").Value = -6573.5

$ws.Range("H31").Value = 2777.0688
$ws.Range("I31").Value = 1178.5264
$ws.Range("J31").Value = 5814.3
$ws.Range("K31").Value = 1178.5264
$ws.Range("L31").Value = 5814.3
$ws.Range("M31").Value = -883.5264
$ws.Range("N31").Value = -6404.3

$ws.Range("H34").Value = 2777.0688
$ws.Range("I34").Value = 1178.5264
$ws.Range("J34").Value = 5814.3
$ws.Range("K34").Value = 1178.5264
$ws.Range("L34").Value = 5814.3
$ws.Range("M34").Value = -976.5264
$ws.Range("N34").Value = -6218.3

$ws.Range("H99").Value = 4160549.8
$ws.Range("I99").Value = 1529875.9
$ws.Range("J99").Value = 5268201.5
$ws.Range("K99").Value = 1529875.9
$ws.Range("L99").Value = 5268201.5
$ws.Range("M99").Value = -1528377.9
$ws.Range("N99").Value = -5271197.5

$ws.Range("H122").Value = 301813.25
$ws.Range("I122").Value = 341636.72
$ws.Range("J122").Value = 3137.25
$ws.Range("K122").Value = 1024910.16
$ws.Range("L122").Value = 9411.75
$ws.Range("M122").Value = -1022460.16
$ws.Range("N122").Value = -14311.75

$ws.Range("H126").Value = 4160549.8
$ws.Range("I126").Value = 1529875.9
$ws.Range("J126").Value = 5268201.5
$ws.Range("K126").Value = 4589627.699999999
$ws.Range("L126").Value = 15804604.5
$ws.Range("M126").Value = -4587157.699999999
$ws.Range("N126").Value = -15809544.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5072
$ws.Range("J80").Value = 4840.625
$ws.Range("L80").Value = 14521.875
$ws.Range("N80").Value = -16393.875

$ws.Range("H83").Value = 5072
$ws.Range("J83").Value = 4840.625
$ws.Range("L83").Value = 43565.625
$ws.Range("N83").Value = -52925.625

$ws.Range("H131").Value = 7932.7144
$ws.Range("I131").Value = 4176.6665
$ws.Range("K131").Value = 12529.9995
$ws.Range("M131").Value = -7489.999500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 82639.28999999999
$ws.Range("I80").Value = 114060
$ws.Range("J80").Value = 4087.5
$ws.Range("K80").Value = 114060
$ws.Range("L80").Value = 4087.5
$ws.Range("M80").Value = -113062
$ws.Range("N80").Value = -6083.5

$ws.Range("H83").Value = 82639.28999999999
$ws.Range("I83").Value = 114060
$ws.Range("J83").Value = 4087.5
$ws.Range("K83").Value = 570300
$ws.Range("L83").Value = 20437.5
$ws.Range("M83").Value = -565308
$ws.Range("N83").Value = -30421.5

$ws.Range("H113").Value = 6072.9
$ws.Range("I113").Value = 3318.5
$ws.Range("J113").Value = 12499.833
$ws.Range("K113").Value = 3318.5
$ws.Range("L113").Value = 12499.833
$ws.Range("M113").Value = -1148.5
$ws.Range("N113").Value = -16839.833

$ws.Range("H122").Value = 3568.6667
$ws.Range("I122").Value = 3554.9524
$ws.Range("J122").Value = 3616.6667
$ws.Range("K122").Value = 10664.8572
$ws.Range("L122").Value = 10850.0001
$ws.Range("M122").Value = -8214.8572
$ws.Range("N122").Value = -15750.0001

$ws.Range("H132").Value = 4476.278
$ws.Range("I132").Value = 4271.846
$ws.Range("J132").Value = 5007.8
$ws.Range("K132").Value = 12815.538
$ws.Range("L132").Value = 15023.4
$ws.Range("M132").Value = -10285.538
$ws.Range("N132").Value = -20083.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8751.588
$ws.Range("I46").Value = 3092
$ws.Range("J46").Value = 9964.357
$ws.Range("K46").Value = 3092
$ws.Range("L46").Value = 9964.357
$ws.Range("M46").Value = -2904
$ws.Range("N46").Value = -10340.357

$ws.Range("H68").Value = 6876.316
$ws.Range("J68").Value = 7176.6665
$ws.Range("L68").Value = 7176.6665
$ws.Range("N68").Value = -8674.666499999999

$ws.Range("H71").Value = 6876.316
$ws.Range("J71").Value = 7176.6665
$ws.Range("L71").Value = 35883.3325
$ws.Range("N71").Value = -43371.3325

$ws.Range("H122").Value = 5364.091
$ws.Range("I122").Value = 4022.1875
$ws.Range("J122").Value = 8942.5
$ws.Range("K122").Value = 12066.5625
$ws.Range("L122").Value = 26827.5
$ws.Range("M122").Value = -9616.5625
$ws.Range("N122").Value = -31727.5

$ws.Range("H132").Value = 3196.8367
$ws.Range("I132").Value = 2478.5405
$ws.Range("J132").Value = 5411.5835
$ws.Range("K132").Value = 7435.6215
$ws.Range("L132").Value = 16234.7505
$ws.Range("M132").Value = -4905.6215
$ws.Range("N132").Value = -21294.7505

$ws.Range("H136").Value = 4685.273
$ws.Range("I136").Value = 2946.5
$ws.Range("J136").Value = 5678.857
$ws.Range("K136").Value = 8839.5
$ws.Range("L136").Value = 17036.571
$ws.Range("M136").Value = -6289.5
$ws.Range("N136").Value = -22136.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4849.2905
$ws.Range("I81").Value = 5913.478
$ws.Range("J81").Value = 1789.75
$ws.Range("K81").Value = 11826.956
$ws.Range("L81").Value = 3579.5
$ws.Range("M81").Value = -10765.956
$ws.Range("N81").Value = -5701.5

$ws.Range("H84").Value = 4849.2905
$ws.Range("I84").Value = 5913.478
$ws.Range("J84").Value = 1789.75
$ws.Range("K84").Value = 59134.78
$ws.Range("L84").Value = 17897.5
$ws.Range("M84").Value = -53830.78
$ws.Range("N84").Value = -28505.5

$ws.Range("H122").Value = 4564.8486
$ws.Range("I122").Value = 1967.4
$ws.Range("J122").Value = 12681.875
$ws.Range("K122").Value = 5902.200000000001
$ws.Range("L122").Value = 38045.625
$ws.Range("M122").Value = -3452.200000000001
$ws.Range("N122").Value = -42945.625

$ws.Range("H126").Value = 1852.9231
$ws.Range("I126").Value = 1852.9231
$ws.Range("K126").Value = 5558.7693
$ws.Range("M126").Value = -3088.7693

$ws.Range("H132").Value = 2442.8262
$ws.Range("I132").Value = 2029.6923
$ws.Range("K132").Value = 6089.0769
$ws.Range("M132").Value = -3559.0769
